$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 21-22: resource/number table ---
$ws.Range("A21").Value = "resource"
$ws.Range("B21").Value = "desert"
$ws.Range("C21").Value = "wheat"
$ws.Range("D21").Value = "whool"
$ws.Range("E21").Value = "clay"
$ws.Range("F21").Value = "ore"
$ws.Range("G21").Value = "wood"

$ws.Range("A22").Value = "number"
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = 3
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 4

# --- Row 27-34: gameboard math helpers ---
$ws.Range("A27").Value = "q"
$ws.Range("B27").Value = 2

$ws.Range("A28").Value = "r"
$ws.Range("B28").Value = 2

$ws.Range("A29").Value = "s"
$ws.Range("F29").Value = "q axis"
$ws.Range("G29").Value = "r pos"

$ws.Range("B30").Value = "L = q*2+1"
$ws.Range("F30").Value = "q axis"
$ws.Range("G30").Value = "r neg"

$ws.Range("B31").Value = "2*(L-1) for r time"
$ws.Range("F31").Value = "r axis"
$ws.Range("G31").Value = "s pos"

$ws.Range("F32").Value = "r axis"
$ws.Range("G32").Value = "s neg"

$ws.Range("F33").Value = "s axis"
$ws.Range("G33").Value = "q pos"

$ws.Range("F34").Value = "s axis"
$ws.Range("G34").Value = "q neg"

# --- Row 38-40: coordinate triples ---
$ws.Range("B38").Value = 1
$ws.Range("C38").Value = -2
$ws.Range("D38").Value = 1

$ws.Range("B39").Value = 1
$ws.Range("C39").Value = -3
$ws.Range("D39").Value = 2

$ws.Range("B40").Value = 2
$ws.Range("C40").Value = -3
$ws.Range("D40").Value = 1

# --- view state tweaks ---
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("E42").Select()
